$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old merged "Note" banner row and drop all old sample data.
$ws.Range("A5:G5").UnMerge()
$ws.Cells.Clear()

# New header row (row 1) - Schedule Description, UserId, Start/End Date, Notes.
$ws.Range("A1").Value = "Schedule Description"
$ws.Range("B1").Value = "UserId"
$ws.Range("C1").Value = "Start Date MM/dd/yyyy"
$ws.Range("D1").Value = "End Date MM/dd/yyyy"
$ws.Range("E1").Value = "Notes"

# Row 2: placeholder text-formatted cells for D:F (matches template's blank data row).
$ws.Range("D2:F2").NumberFormat = "@"

# Approximate the new bestFit-calculated column widths (closest reachable via ColumnWidth).
$ws.Columns.Item(2).ColumnWidth = 6
$ws.Columns.Item(3).ColumnWidth = 21.5
$ws.Columns.Item(4).ColumnWidth = 20.666667
$ws.Columns.Item(5).ColumnWidth = 5.5
$ws.Columns.Item(6).ColumnWidth = 5.5

# Update the remembered selection on the sheet.
$ws.Range("E8").Select()
